$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 50
$ws.Range("G5").Value = 7
$ws.Range("H5").Value = 1
$ws.Range("G6").Value = 8
$ws.Range("H6").Value = 1

$ws.Range("G3").Select()
